$d = $word.ActiveDocument

$d.Content.Find.Execute("962÷6=160, 2", $true, $false, $false, $false, $false, $true, 1, $false, "775÷2=387, 1", 2)
$d.Content.Find.Execute("638÷9=70, 8", $true, $false, $false, $false, $false, $true, 1, $false, "940÷3=313, 1", 2)
$d.Content.Find.Execute("352÷4=88, 0", $true, $false, $false, $false, $false, $true, 1, $false, "100÷7=14, 2", 2)
$d.Content.Find.Execute("318÷5=63, 3", $true, $false, $false, $false, $false, $true, 1, $false, "200÷5=40, 0", 2)
$d.Content.Find.Execute("682÷2=341, 0", $true, $false, $false, $false, $false, $true, 1, $false, "306÷2=153, 0", 2)
$d.Content.Find.Execute("560÷6=93, 2", $true, $false, $false, $false, $false, $true, 1, $false, "551÷7=78, 5", 2)
$d.Content.Find.Execute("454÷5=90, 4", $true, $false, $false, $false, $false, $true, 1, $false, "466÷2=233, 0", 2)
$d.Content.Find.Execute("935÷2=467, 1", $true, $false, $false, $false, $false, $true, 1, $false, "613÷2=306, 1", 2)
$d.Content.Find.Execute("243÷6=40, 3", $true, $false, $false, $false, $false, $true, 1, $false, "169÷5=33, 4", 2)
$d.Content.Find.Execute("232÷9=25, 7", $true, $false, $false, $false, $false, $true, 1, $false, "904÷3=301, 1", 2)
$d.Content.Find.Execute("443÷8=55, 3", $true, $false, $false, $false, $false, $true, 1, $false, "176÷8=22, 0", 2)
$d.Content.Find.Execute("188÷9=20, 8", $true, $false, $false, $false, $false, $true, 1, $false, "654÷2=327, 0", 2)
$d.Content.Find.Execute("995÷2=497, 1", $true, $false, $false, $false, $false, $true, 1, $false, "868÷2=434, 0", 2)
$d.Content.Find.Execute("346÷2=173, 0", $true, $false, $false, $false, $false, $true, 1, $false, "195÷7=27, 6", 2)
$d.Content.Find.Execute("252÷3=84, 0", $true, $false, $false, $false, $false, $true, 1, $false, "824÷9=91, 5", 2)
$d.Content.Find.Execute("151÷4=37, 3", $true, $false, $false, $false, $false, $true, 1, $false, "463÷2=231, 1", 2)
$d.Content.Find.Execute("217÷8=27, 1", $true, $false, $false, $false, $false, $true, 1, $false, "920÷8=115, 0", 2)
$d.Content.Find.Execute("154÷9=17, 1", $true, $false, $false, $false, $false, $true, 1, $false, "755÷9=83, 8", 2)
$d.Content.Find.Execute("920÷5=184, 0", $true, $false, $false, $false, $false, $true, 1, $false, "698÷3=232, 2", 2)
$d.Content.Find.Execute("969÷3=323, 0", $true, $false, $false, $false, $false, $true, 1, $false, "385÷8=48, 1", 2)
$d.Content.Find.Execute("485÷7=69, 2", $true, $false, $false, $false, $false, $true, 1, $false, "435÷6=72, 3", 2)
$d.Content.Find.Execute("835÷5=167, 0", $true, $false, $false, $false, $false, $true, 1, $false, "115÷2=57, 1", 2)
$d.Content.Find.Execute("161÷4=40, 1", $true, $false, $false, $false, $false, $true, 1, $false, "965÷8=120, 5", 2)
$d.Content.Find.Execute("971÷9=107, 8", $true, $false, $false, $false, $false, $true, 1, $false, "871÷7=124, 3", 2)
$d.Content.Find.Execute("345÷4=86, 1", $true, $false, $false, $false, $false, $true, 1, $false, "804÷6=134, 0", 2)
